# ex07 p3 - both kernels work - time to get benchmarking
#
# Adds a new worksheet "P3" after the existing "P2" sheet, populates it
# with a small doubling series (A3=30, then A4:A10 = previous*2), sizes
# column C, and leaves it as the active/selected sheet - matching the
# target workbook state.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet (P2), then rename it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "P3"

# Seed value and doubling formulas down column A.
$ws3.Range("A3").Value = 30
$ws3.Range("A4").Formula = "=A3*2"
$ws3.Range("A5").Formula = "=A4*2"
$ws3.Range("A6").Formula = "=A5*2"
$ws3.Range("A7").Formula = "=A6*2"
$ws3.Range("A8").Formula = "=A7*2"
$ws3.Range("A9").Formula = "=A8*2"
$ws3.Range("A10").Formula = "=A9*2"

# Match the stored column width of 10 for column C.
$ws3.Columns.Item(3).ColumnWidth = 9.14

# Leave the selection/active cell on A11, with P3 as the active tab.
$ws3.Range("A11").Select() | Out-Null
